$d = $word.ActiveDocument

# Add a new paragraph at the end of the document with "Thank you" text
$end = $d.Content.End - 1
$range = $d.Range($end, $end)
$range.InsertParagraphAfter()
$range = $d.Range($end, $end)
